$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 39.90000152587891
$ws.Range("E2").Value = 38.70999908447266
$ws.Range("F2").Value = 41.81499862670898
$ws.Range("G2").Value = 37.2400016784668
$ws.Range("H2").Value = 52693110
$ws.Range("I2").Value = "HUBS"

$ws.Range("D3").Value = 49.97999954223633
$ws.Range("E3").Value = 53.95000076293945
$ws.Range("F3").Value = 55.06999969482422
$ws.Range("G3").Value = 46.70999908447266
$ws.Range("H3").Value = 52693110
$ws.Range("I3").Value = "HUBS"

$ws.Range("D4").Value = 46.31999969482422
$ws.Range("E4").Value = 51.88000106811523
$ws.Range("F4").Value = 52.9900016784668
$ws.Range("G4").Value = 44.36000061035156
$ws.Range("H4").Value = 52693110
$ws.Range("I4").Value = "HUBS"

$ws.Range("D5").Value = 55.43999862670898
$ws.Range("E5").Value = 40.59000015258789
$ws.Range("F5").Value = 55.97999954223633
$ws.Range("G5").Value = 38.04999923706055
$ws.Range("H5").Value = 52693110
$ws.Range("I5").Value = "HUBS"

$ws.Range("D6").Value = 43.29999923706055
$ws.Range("E6").Value = 44.29000091552734
$ws.Range("F6").Value = 45.65999984741211
$ws.Range("G6").Value = 39.0099983215332
$ws.Range("H6").Value = 52693110
$ws.Range("I6").Value = "HUBS"

$ws.Range("D7").Value = 43.36000061035156
$ws.Range("E7").Value = 54.59000015258789
$ws.Range("F7").Value = 54.59999847412109
$ws.Range("G7").Value = 41.70999908447266
$ws.Range("H7").Value = 52693110
$ws.Range("I7").Value = "HUBS"

$ws.Range("D8").Value = 57
$ws.Range("E8").Value = 52.45000076293945
$ws.Range("F8").Value = 59
$ws.Range("G8").Value = 51.25
$ws.Range("H8").Value = 52693110
$ws.Range("I8").Value = "HUBS"

$ws.Range("D9").Value = 47.29999923706055
$ws.Range("E9").Value = 51.29999923706055
$ws.Range("F9").Value = 53.5
$ws.Range("G9").Value = 46.79999923706055
$ws.Range("H9").Value = 52693110
$ws.Range("I9").Value = "HUBS"

$ws.Range("D10").Value = 60.54999923706055
$ws.Range("E10").Value = 67.05000305175781
$ws.Range("F10").Value = 67.30000305175781
$ws.Range("G10").Value = 58.45000076293945
$ws.Range("H10").Value = 52693110
$ws.Range("I10").Value = "HUBS"

$ws.Range("D11").Value = 66.30000305175781
$ws.Range("E11").Value = 72.3499984741211
$ws.Range("F11").Value = 76.0999984741211
$ws.Range("G11").Value = 63
$ws.Range("H11").Value = 52693110
$ws.Range("I11").Value = "HUBS"

$ws.Range("D12").Value = 84.1500015258789
$ws.Range("E12").Value = 86.55000305175781
$ws.Range("F12").Value = 87.25
$ws.Range("G12").Value = 81.44999694824219
$ws.Range("H12").Value = 52693110
$ws.Range("I12").Value = "HUBS"

$ws.Range("D13").Value = 88.4000015258789
$ws.Range("E13").Value = 97.0500030517578
$ws.Range("F13").Value = 102.3000030517578
$ws.Range("G13").Value = 87.5999984741211
$ws.Range("H13").Value = 52693110
$ws.Range("I13").Value = "HUBS"

$ws.Range("D14").Value = 107.6500015258789
$ws.Range("E14").Value = 105.9000015258789
$ws.Range("F14").Value = 117.5999984741211
$ws.Range("G14").Value = 101.4499969482422
$ws.Range("H14").Value = 52693110
$ws.Range("I14").Value = "HUBS"

$ws.Range("D15").Value = 124.9499969482422
$ws.Range("E15").Value = 124.0999984741211
$ws.Range("F15").Value = 136.75
$ws.Range("G15").Value = 122
$ws.Range("H15").Value = 52693110
$ws.Range("I15").Value = "HUBS"

$ws.Range("D16").Value = 151.9199981689453
$ws.Range("E16").Value = 135.6499938964844
$ws.Range("F16").Value = 153.5399932861328
$ws.Range("G16").Value = 121.129997253418
$ws.Range("H16").Value = 52693110
$ws.Range("I16").Value = "HUBS"

$ws.Range("D17").Value = 122.5800018310547
$ws.Range("E17").Value = 158.3099975585938
$ws.Range("F17").Value = 158.6600036621094
$ws.Range("G17").Value = 117.7699966430664
$ws.Range("H17").Value = 52693110
$ws.Range("I17").Value = "HUBS"

$ws.Range("D18").Value = 166.25
$ws.Range("E18").Value = 184.4900054931641
$ws.Range("F18").Value = 186.9649963378907
$ws.Range("G18").Value = 157.9900054931641
$ws.Range("H18").Value = 52693110
$ws.Range("I18").Value = "HUBS"

$ws.Range("D19").Value = 175.0599975585938
$ws.Range("E19").Value = 178.7200012207031
$ws.Range("F19").Value = 186.5
$ws.Range("G19").Value = 171.4799957275391
$ws.Range("H19").Value = 52693110
$ws.Range("I19").Value = "HUBS"

$ws.Range("D20").Value = 152.2200012207031
$ws.Range("E20").Value = 155.1000061035156
$ws.Range("F20").Value = 168.4299926757812
$ws.Range("G20").Value = 145.5950012207031
$ws.Range("H20").Value = 52693110
$ws.Range("I20").Value = "HUBS"

$ws.Range("D21").Value = 159.5099945068359
$ws.Range("E21").Value = 180.9400024414062
$ws.Range("F21").Value = 188.8500061035156
$ws.Range("G21").Value = 159.1000061035156
$ws.Range("H21").Value = 52693110
$ws.Range("I21").Value = "HUBS"

$ws.Range("D22").Value = 127.5
$ws.Range("E22").Value = 168.6300048828125
$ws.Range("F22").Value = 173.2949981689453
$ws.Range("G22").Value = 110.5749969482422
$ws.Range("H22").Value = 52693110
$ws.Range("I22").Value = "HUBS"

$ws.Range("D23").Value = 223.979995727539
$ws.Range("E23").Value = 234.6100006103516
$ws.Range("F23").Value = 243.5800018310547
$ws.Range("G23").Value = 205.0700073242188
$ws.Range("H23").Value = 52693110
$ws.Range("I23").Value = "HUBS"

$ws.Range("D24").Value = 296.1900024414062
$ws.Range("E24").Value = 290.0700073242188
$ws.Range("F24").Value = 330.7160034179688
$ws.Range("G24").Value = 281.5700073242188
$ws.Range("H24").Value = 52693110
$ws.Range("I24").Value = "HUBS"

$ws.Range("D25").Value = 396.4400024414063
$ws.Range("E25").Value = 372.2000122070313
$ws.Range("F25").Value = 414.6740112304688
$ws.Range("G25").Value = 347.7799987792969
$ws.Range("H25").Value = 52693110
$ws.Range("I25").Value = "HUBS"

$ws.Range("D26").Value = 465.1700134277344
$ws.Range("E26").Value = 526.4500122070312
$ws.Range("F26").Value = 574.8300170898438
$ws.Range("G26").Value = 463.3500061035156
$ws.Range("H26").Value = 52693110
$ws.Range("I26").Value = "HUBS"

$ws.Range("D27").Value = 578.8200073242188
$ws.Range("E27").Value = 596.02001953125
$ws.Range("F27").Value = 616.4500122070312
$ws.Range("G27").Value = 543.1699829101562
$ws.Range("H27").Value = 52693110
$ws.Range("I27").Value = "HUBS"

$ws.Range("D28").Value = 675.3599853515625
$ws.Range("E28").Value = 810.22998046875
$ws.Range("F28").Value = 839.7750244140625
$ws.Range("G28").Value = 622.0399780273438
$ws.Range("H28").Value = 52693110
$ws.Range("I28").Value = "HUBS"

$ws.Range("D29").Value = 659.239990234375
$ws.Range("E29").Value = 488.7999877929688
$ws.Range("F29").Value = 659.239990234375
$ws.Range("G29").Value = 403
$ws.Range("H29").Value = 52693110
$ws.Range("I29").Value = "HUBS"

$ws.Range("D30").Value = 477.239990234375
$ws.Range("E30").Value = 379.4299926757813
$ws.Range("F30").Value = 514.3400268554688
$ws.Range("G30").Value = 367.1050109863281
$ws.Range("H30").Value = 52693110
$ws.Range("I30").Value = "HUBS"

$ws.Range("D31").Value = 301.6400146484375
$ws.Range("E31").Value = 308
$ws.Range("F31").Value = 335.75
$ws.Range("G31").Value = 257.2099914550781
$ws.Range("H31").Value = 52693110
$ws.Range("I31").Value = "HUBS"

$ws.Range("D32").Value = 271.0499877929688
$ws.Range("E32").Value = 296.5599975585937
$ws.Range("F32").Value = 305.8800048828125
$ws.Range("G32").Value = 245.0299987792969
$ws.Range("H32").Value = 52693110
$ws.Range("I32").Value = "HUBS"

$ws.Range("D33").Value = 297.7699890136719
$ws.Range("E33").Value = 347.010009765625
$ws.Range("F33").Value = 362.5249938964844
$ws.Range("G33").Value = 267.989990234375
$ws.Range("H33").Value = 52693110
$ws.Range("I33").Value = "HUBS"

$ws.Range("D34").Value = 423.3399963378906
$ws.Range("E34").Value = 420.9500122070313
$ws.Range("F34").Value = 429.5199890136719
$ws.Range("G34").Value = 390.2200012207031
$ws.Range("H34").Value = 52693110
$ws.Range("I34").Value = "HUBS"

$ws.Range("D35").Value = 529.6199951171875
$ws.Range("E35").Value = 580.5499877929688
$ws.Range("F35").Value = 581.4000244140625
$ws.Range("G35").Value = 497.0280151367188
$ws.Range("H35").Value = 52693110
$ws.Range("I35").Value = "HUBS"

$ws.Range("D36").Value = 492.1199951171875
$ws.Range("E36").Value = 423.7699890136719
$ws.Range("F36").Value = 505.2900085449219
$ws.Range("G36").Value = 407.2300109863281
$ws.Range("H36").Value = 52693110
$ws.Range("I36").Value = "HUBS"

$ws.Range("D37").Value = 572.9400024414062
$ws.Range("E37").Value = 611
$ws.Range("F37").Value = 635.739990234375
$ws.Range("G37").Value = 527
$ws.Range("H37").Value = 52693110
$ws.Range("I37").Value = "HUBS"

$ws.Range("D38").Value = 624.0700073242188
$ws.Range("E38").Value = 604.8699951171875
$ws.Range("F38").Value = 693.8499755859375
$ws.Range("G38").Value = 604.5700073242188
$ws.Range("H38").Value = 52693110
$ws.Range("I38").Value = "HUBS"

$ws.Range("D39").Value = 585
$ws.Range("E39").Value = 497.0299987792969
$ws.Range("F39").Value = 596.25
$ws.Range("G39").Value = 455
$ws.Range("H39").Value = 52693110
$ws.Range("I39").Value = "HUBS"

$ws.Range("D40").Value = 532.2899780273438
$ws.Range("E40").Value = 554.7899780273438
$ws.Range("F40").Value = 570.6500244140625
$ws.Range("G40").Value = 510.1199951171875
$ws.Range("H40").Value = 52693110
$ws.Range("I40").Value = "HUBS"

$ws.Range("D41").Value = 704.97998046875
$ws.Range("E41").Value = 779.530029296875
$ws.Range("F41").Value = 811.7899780273438
$ws.Range("G41").Value = 685.6749877929688
$ws.Range("H41").Value = 52693110
$ws.Range("I41").Value = "HUBS"

$ws.Range("D42").Value = 573.719970703125
$ws.Range("E42").Value = 611.5
$ws.Range("F42").Value = 622.3800048828125
$ws.Range("G42").Value = 472.2780151367188
$ws.Range("H42").Value = 52693110
$ws.Range("I42").Value = "HUBS"

$ws.Range("D43").Value = 555.02001953125
$ws.Range("E43").Value = 519.6500244140625
$ws.Range("F43").Value = 568.1599731445312
$ws.Range("G43").Value = 516.6300048828125
$ws.Range("H43").Value = 52693110
$ws.Range("I43").Value = "HUBS"

Write-Host "Edit complete"